$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.378.40'
$ws.Range('E2').Value = '  -2.65%  '
$ws.Range('D3').Value = '2.420.45'
$ws.Range('E3').Value = '  -3.54%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '511.64'
$ws.Range('E5').Value = '  -3.99%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '128.69'
$ws.Range('E6').Value = '  -4.85%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.549'
$ws.Range('E8').Value = '  -3.11%  '
$ws.Range('D9').Value = '2.428.76'
$ws.Range('E9').Value = '  -3.37%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.156'
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.20'
$ws.Range('E12').Value = '  -3.89%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.332'
$ws.Range('E13').Value = '  -3.82%  '
$ws.Range('D14').Value = '2.850.94'
$ws.Range('E14').Value = '  -3.48%  '
$ws.Range('D15').Value = '57.316.80'
$ws.Range('E15').Value = '  -2.59%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.46'
$ws.Range('E16').Value = '  -5.68%  '
$ws.Range('E17').Value = '  -4.63%  '
$ws.Range('D18').Value = '2.424.73'
$ws.Range('E18').Value = '  -3.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.36'
$ws.Range('E19').Value = '  -5.95%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '314.35'
$ws.Range('E20').Value = '  -2.45%  '
$ws.Range('E21').Value = '  -4.16%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.65'
$ws.Range('E23').Value = '  -4.98%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.52'
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.402'
$ws.Range('E25').Value = '  -4.32%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -2.61%  '
$ws.Range('E28').Value = '  -4.87%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '169.02'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '0.0₃0717'
$ws.Range('E30').Value = '  -5.95%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.17'
$ws.Range('E31').Value = '  -4.80%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.65'
$ws.Range('E32').Value = '  -5.17%  '
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.68'
$ws.Range('E36').Value = '  -4.02%  '
$ws.Range('E37').Value = '  -7.29%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.85'
$ws.Range('E38').Value = '  -4.60%  '
$ws.Range('E39').Value = '  -2.08%  '
$ws.Range('E40').Value = '  -5.55%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.764'
$ws.Range('E41').Value = '  -4.19%  '
$ws.Range('E42').Value = '  -5.76%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '264.43'
$ws.Range('E43').Value = '  -5.93%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.86'
$ws.Range('E44').Value = '  -2.93%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.582'
$ws.Range('E45').Value = '  -3.68%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '121.17'
$ws.Range('E46').Value = '  -6.55%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0901'
$ws.Range('E47').Value = '  -2.67%  '
$ws.Range('E48').Value = '  -3.81%  '
$ws.Range('E49').Value = '  -3.76%  '
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('D51').Value = '1.690.32'
$ws.Range('E51').Value = '  -3.82%  '
